# adding new practice scenarios
# Turn the A2 "username" value into a mailto: hyperlink pointing at a new
# practice email address, matching the target XML diff:
#  - new shared string "automationuser1982@gmail.com"
#  - new Hyperlink-style font/xf/cellStyle in styles.xml
#  - A2 becomes a shared-string cell styled with the Hyperlink cell style
#  - a <hyperlinks> entry on the worksheet referencing the new relationship
#  - A2 becomes the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:automationuser1982@gmail.com")
$ws.Range("A2").Value = "automationuser1982@gmail.com"
[void]$ws.Range("A2").Select()
